$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.32"
$ws.Range("E2").Value = "'1.11%"
$ws.Range("D3").Value = "'27.11"
$ws.Range("E3").Value = "'0.84%"
$ws.Range("D4").Value = "'4.704"
$ws.Range("E4").Value = "'1.23%"
$ws.Range("D5").Value = "'0.06185"
$ws.Range("E5").Value = "'3.41%"
$ws.Range("D6").Value = "'6.687"
$ws.Range("E6").Value = "'0.63%"
$ws.Range("E7").Value = "'-0.71%"
$ws.Range("D8").Value = "'0.9151"
$ws.Range("E8").Value = "'-0.57%"
$ws.Range("D9").Value = "'0.1407"
$ws.Range("E9").Value = "'1.46%"
$ws.Range("D10").Value = "'0.04668"
$ws.Range("E10").Value = "'-3.85%"
$ws.Range("D11").Value = "'0.07087"
$ws.Range("E11").Value = "'0.91%"
$ws.Range("D12").Value = "'0.03150"
$ws.Range("E12").Value = "'3.26%"
$ws.Range("D13").Value = "'0.09045"
$ws.Range("E13").Value = "'-0.80%"
$ws.Range("D14").Value = "'0.001533"
$ws.Range("E14").Value = "'-0.24%"
$ws.Range("D15").Value = "'0.0006169"
$ws.Range("E15").Value = "'1.64%"
$ws.Range("D16").Value = "'0.006129"
$ws.Range("E16").Value = "'-1.06%"
$ws.Range("E17").Value = "'0.31%"
$ws.Range("D18").Value = "'3.177"
$ws.Range("E18").Value = "'0.93%"
$ws.Range("E19").Value = "'-1.72%"
$ws.Range("E20").Value = "'-0.95%"
$ws.Range("E21").Value = "'0.85%"
$ws.Range("D22").Value = "'4.085"
$ws.Range("E22").Value = "'0.93%"
$ws.Range("D23").Value = "'0.04224"
$ws.Range("E23").Value = "'0.24%"
$ws.Range("E24").Value = "'-0.14%"
$ws.Range("E25").Value = "'-5.55%"
$ws.Range("E26").Value = "'0.06%"
$ws.Range("D27").Value = "'0.0001578"
$ws.Range("E27").Value = "'-7.82%"
$ws.Range("E40").Value = "'1.73%"
$ws.Range("D41").Value = "'0.1112"
$ws.Range("E41").Value = "'0.04%"
$ws.Range("D42").Value = "'0.004099"
$ws.Range("E42").Value = "'8.77%"
$ws.Range("E43").Value = "'7.41%"
$ws.Range("D44").Value = "'0.002184"
$ws.Range("E44").Value = "'-10.09%"
$ws.Range("D45").Value = "'0.00005159"
$ws.Range("E45").Value = "'0.92%"
$ws.Range("E46").Value = "'0.06%"
$ws.Range("D48").Value = "'0.1677"
$ws.Range("E48").Value = "'56.37%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.06%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.06%"
